$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: columns A-T (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then 16 numeric metric columns). This replaces the
# existing rows 2-7 content and appends new rows 8-10, expanding the
# sending/target cluster combinations to the full 3x3 grid (ECs/FAPs/sCs).
$data = @(
  @("ECs","Sema4d","Erbb2","ECs",3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,165.8425359564111,1492.5828236077,0.2252471739205623,0.2252471739205623),
  @("ECs","Sema4d","Erbb2","FAPs",3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,322.4226075905267,2901.80346831474,0.4379140776462364,0.4379140776462364),
  @("ECs","Sema4d","Erbb2","sCs",3,1,76.08252666666667,228.24758,0.95878149807566,0.95878149807566,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,217.6560553802067,1958.90449842186,0.2956202465088614,0.2956202465088613),
  @("FAPs","Sema4d","Erbb2","ECs",3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,2.037521947471667,18.337697527245,0.002767360362782558,0.002767360362782558),
  @("FAPs","Sema4d","Erbb2","FAPs",3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,3.961246344541,35.651217100869,0.005380161001309901,0.005380161001309901),
  @("FAPs","Sema4d","Erbb2","sCs",3,1,0.9347409999999999,2.804223,0.01177947704364805,0.01177947704364805,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,2.674096770649,24.066870935841,0.003631955679555589,0.003631955679555589),
  @("sCs","Sema4d","Erbb2","ECs",3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,5.092132620515557,45.82919358464,0.006916129661097937,0.006916129661097936),
  @("sCs","Sema4d","Erbb2","FAPs",3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,9.899864761685334,89.09878285516801,0.01344598686281231,0.01344598686281231),
  @("sCs","Sema4d","Erbb2","sCs",3,1,2.336085333333334,7.008256,0.02943902488069198,0.02943902488069198,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,6.683047224661335,60.14742502195201,0.009076908356781731,0.009076908356781731)
)

$r = 2
foreach ($row in $data) {
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
  $r++
}
